$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for Ids 105030033 / 105030138 / 105030139 (currently on
# sheet rows 3, 4 and 5 respectively) get cyclically rotated:
#   new row 3 <- old row 4   (Id 105030138)
#   new row 4 <- old row 5   (Id 105030139)
#   new row 5 <- old row 3   (Id 105030033)
#
# Use Range.Copy (cell-to-cell copy) rather than reading/writing .Value,
# so text-typed cells (e.g. the Startdatum/Slutdatum date strings) keep
# their original text type instead of being reinterpreted as dates, and
# no new/extra cell styles get introduced.
#
# A scratch row far outside the used range (A1:AY6) holds old row 3
# while rows 4 and 5 are shifted up.

$ws.Range("A3:AY3").Copy($ws.Range("A100:AY100"))   # scratch <- old row 3
$ws.Range("A4:AY4").Copy($ws.Range("A3:AY3"))        # row 3   <- old row 4
$ws.Range("A5:AY5").Copy($ws.Range("A4:AY4"))        # row 4   <- old row 5
$ws.Range("A100:AY100").Copy($ws.Range("A5:AY5"))    # row 5   <- scratch (old row 3)

# Clean up the scratch row.
$ws.Range("A100:AY100").Clear()
